# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp string in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 10:49"

# --- Swap Groenlandia (row 209) and Islas Malvinas (row 210) ---
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(210, 1).Value = "Groenlandia"

# --- Update numeric stats for affected countries ---
# Row 7: Rusia
$ws.Cells.Item(7, 2).Value = 733699
$ws.Cells.Item(7, 3).Value = 6537
$ws.Cells.Item(7, 4).Value = 504021
$ws.Cells.Item(7, 5).Value = 218239
$ws.Cells.Item(7, 7).Value = 104
$ws.Cells.Item(7, 8).Value = 11439

# Row 47: Polonia
$ws.Cells.Item(47, 2).Value = 38190
$ws.Cells.Item(47, 3).Value = 299
$ws.Cells.Item(47, 4).Value = 27515
$ws.Cells.Item(47, 5).Value = 9099
$ws.Cells.Item(47, 7).Value = 5
$ws.Cells.Item(47, 8).Value = 1576

# Row 48: Afganistan
$ws.Cells.Item(48, 2).Value = 34455
$ws.Cells.Item(48, 3).Value = 4
$ws.Cells.Item(48, 4).Value = 21254
$ws.Cells.Item(48, 5).Value = 12189
$ws.Cells.Item(48, 7).Value = 2
$ws.Cells.Item(48, 8).Value = 1012

# Row 60: Moldavia
$ws.Cells.Item(60, 4).Value = 12793
$ws.Cells.Item(60, 5).Value = 5947

# Row 62: Austria
$ws.Cells.Item(62, 2).Value = 18948
$ws.Cells.Item(62, 3).Value = 51
$ws.Cells.Item(62, 4).Value = 17000
$ws.Cells.Item(62, 5).Value = 1240

# Row 118: Eslovaquia
$ws.Cells.Item(118, 2).Value = 1902
$ws.Cells.Item(118, 3).Value = 1
$ws.Cells.Item(118, 5).Value = 381

# Row 140: Uganda
$ws.Cells.Item(140, 2).Value = 1029
$ws.Cells.Item(140, 3).Value = 4
$ws.Cells.Item(140, 4).Value = 977
$ws.Cells.Item(140, 5).Value = 52

# Row 190: Islas Turcas y Caicos
$ws.Cells.Item(190, 2).Value = 72
$ws.Cells.Item(190, 3).Value = 1
$ws.Cells.Item(190, 5).Value = 59
